$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Identificadores de CU (orden en que se agregaron las cadenas compartidas)
$ws.Range("B27").Value = "CU - 23"
$ws.Range("B28").Value = "CU - 24"

# Alias / título corto
$ws.Range("D27").Value = "Consultar profesores"
$ws.Range("D28").Value = "Consultar clientes"

# Descripciones
$ws.Range("C27").Value = "El director puede  consultar todos los profesores en la institución."
$ws.Range("C28").Value = "El director puede consultar todos los clientes registrados."

# Estado
$ws.Range("E27").Value = "vacio"
$ws.Range("E28").Value = "vacio"

# Esfuerzo (hrs)
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0

# Incremento
$ws.Range("G27").Value = 0
$ws.Range("G28").Value = 0

# Prioridad
$ws.Range("H27").Value = 1
$ws.Range("H28").Value = 1

$ws.Range("C28").Select()
